$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2022-09-02 21:00:58"

for ($row = 2; $row -le 64; $row++) {
    $ws.Range("O$row").Value = $newTimestamp
}
